$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T55")

# Row 2
$ws.Range("B2").Value = -0.08279561087709399
$ws.Range("C2").Value = 2.111406769062893
$ws.Range("D2").Value = 17.09390349523057
$ws.Range("E2").Value = 4.134477415010338
$ws.Range("F2").Value = 4.230923936047997
$ws.Range("G2").Value = 22

# Row 3
$ws.Range("B3").Value = 0.9006556730100386
$ws.Range("C3").Value = 2.200387812214267
$ws.Range("D3").Value = 19.30760674930202
$ws.Range("E3").Value = 4.394042187929244
$ws.Range("F3").Value = 4.397405229137807
$ws.Range("G3").Value = 23

# Row 4
$ws.Range("B4").Value = -0.1540770173485546
$ws.Range("C4").Value = 0.9992936523105611
$ws.Range("D4").Value = 3.826009492865405
$ws.Range("E4").Value = 1.956018786429569
$ws.Range("F4").Value = 2.000597217567105
$ws.Range("G4").Value = 20
